# Weekly update: a new week of "Apio" (celery) price data was added at the
# top of the data block (row 150). All existing records from row 150
# through row 189 shift down by one row (to 151-190), and the new row 150
# is populated with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 150, pushing rows 150:189 down to 151:190.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row 150 with this week's data.
$ws.Range("A150").Value = 5
$ws.Range("B150").Value = "Macroferia Regional de Talca"
$ws.Range("C150").Value = "Maule"
$ws.Range("D150").Value = 44736
$ws.Range("E150").Value = 7
$ws.Range("F150").Value = 100112017
$ws.Range("G150").Value = "Apio"
$ws.Range("H150").Value = "Americana (o)"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 700
$ws.Range("K150").Value = 7000
$ws.Range("L150").Value = 7000
$ws.Range("M150").Value = 7000
$ws.Range("N150").Value = "`$/docena de matas"
$ws.Range("O150").Value = "Provincia del Elquí"
$ws.Range("P150").Value = 1167
$ws.Range("Q150").Value = 6
$ws.Range("R150").Value = "Hortaliza"
